# Spring 23 week 14 inputs - update SL matchup average values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.38
$ws.Range("E2").Value = 1.3

$ws.Range("B3").Value = 1.46
$ws.Range("D3").Value = 1.44

$ws.Range("E4").Value = 1.25

$ws.Range("B5").Value = 1.48
$ws.Range("G5").Value = 0.68

$ws.Range("E6").Value = 1.32

$ws.Range("E7").Value = 1.94
